# Update "想去人数" (number of people interested) figures that changed
# between the previous data scrape and the new one generated at commit 456a3b4.
#
# Sheet "展览" (exhibitions):
#   F5: 4024 -> 4034
#   F6:   35 ->   36
#   F7:  445 ->  446
#
# Sheet "全部类型" (all types - aggregate of every category):
#   F5: 4024 -> 4034
#   F8:   35 ->   36
#   F9:  445 ->  446

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 4034
$wsExhibit.Range("F6").Value = 36
$wsExhibit.Range("F7").Value = 446

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 4034
$wsAll.Range("F8").Value = 36
$wsAll.Range("F9").Value = 446
